# Update gh-pages to output generated at 456a3b4
#
# Refreshes the "想去人数" (column F) counters across the three source
# sheets ("展览", "演出") and the aggregated "全部类型" sheet that mirrors
# every event. Two events that had sold out ("不可售" in column G) are
# back on sale, so their "最低票价" (column G) is restored to a numeric
# price as well.

$wb = $excel.ActiveWorkbook

$sheetExpo = $wb.Worksheets.Item(1)   # 展览
$sheetShow = $wb.Worksheets.Item(2)   # 演出
$sheetAll  = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1: 展览 -----------------------------------------------------
$sheetExpo.Range("F3").Value  = 26
$sheetExpo.Range("F5").Value  = 15478
$sheetExpo.Range("F8").Value  = 695
$sheetExpo.Range("F9").Value  = 15372
$sheetExpo.Range("F11").Value = 8971
$sheetExpo.Range("F12").Value = 367
$sheetExpo.Range("F20").Value = 43
$sheetExpo.Range("F21").Value = 541
$sheetExpo.Range("F25").Value = 1104
$sheetExpo.Range("F33").Value = 38

# Row 35 went from sold out ("不可售") back on sale at 78
$sheetExpo.Range("F35").Value = 306
$sheetExpo.Range("G35").Value = 78

$sheetExpo.Range("F36").Value = 445
$sheetExpo.Range("F38").Value = 5496

# --- Sheet 2: 演出 -----------------------------------------------------
$sheetShow.Range("F2").Value = 66

# --- Sheet 4: 全部类型 --------------------------------------------------
$sheetAll.Range("F3").Value  = 26
$sheetAll.Range("F5").Value  = 15478
$sheetAll.Range("F8").Value  = 695
$sheetAll.Range("F9").Value  = 15372
$sheetAll.Range("F11").Value = 8971
$sheetAll.Range("F12").Value = 367
$sheetAll.Range("F20").Value = 43
$sheetAll.Range("F21").Value = 541
$sheetAll.Range("F25").Value = 1104
$sheetAll.Range("F31").Value = 66
$sheetAll.Range("F35").Value = 38

# Row 37 went from sold out ("不可售") back on sale at 78
$sheetAll.Range("F37").Value = 306
$sheetAll.Range("G37").Value = 78

$sheetAll.Range("F38").Value = 445
$sheetAll.Range("F40").Value = 5496
